# Add "x" markers into column A for the module rows that didn't have one yet
# (rows 32, 35, 36, 38 on the "Main" sheet - the rows with the grey/shaded
# B/C styling), then leave the selection where the user had it when they
# saved (cell A34, scrolled so row 28 is the top visible row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Activate()

$ws.Range("A32").Value = "x"
$ws.Range("A35").Value = "x"
$ws.Range("A36").Value = "x"
$ws.Range("A38").Value = "x"

# Best-effort: scroll the window so row 28 is the top-left visible cell,
# matching the saved view. (May be a no-op in headless/sandbox hosts that
# don't track window scroll position, but is harmless either way.)
try {
    $excel.ActiveWindow.ScrollRow = 28
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$ws.Range("A34").Select()
